$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TER")

# Row 4 - Inventory
$ws.Range("B4").Value = 262000000.0
$ws.Range("C4").Value = 222000000.0
$ws.Range("D4").Value = 191000000.0
$ws.Range("E4").Value = 206000000.0
$ws.Range("F4").Value = 183000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 177000000.0
$ws.Range("C15").Value = 134000000.0
$ws.Range("D15").Value = 150000000.0
$ws.Range("E15").Value = 184000000.0
$ws.Range("F15").Value = 130000000.0

# Row 26 - Long Term Tax Liability (Deferred)
$ws.Range("B26").Value = -83000000.0
$ws.Range("C26").Value = -77000000.0
$ws.Range("D26").Value = -67000000.0
$ws.Range("E26").Value = -67000000.0
$ws.Range("F26").Value = -63000000.0

# Row 38 - Net Debt
$ws.Range("G38").Value = -497540000.0

# Row 39 - Total Debt
$ws.Range("G39").Value = 413687000.0
